$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("J7").Value = "upper upper"
$ws.Range("K11").Value = 123
Write-Output "done"
